# Generate Report for Handback
# The b9d5d801-102b-48d6-9175-1fbd789d7af6 row has finished its handback
# cycle: flip its "Ready for handoff" status to "Handed back: in sync with
# en-US" on every sheet, and stamp the new handback timestamps on the
# per-locale sheets.

$wb = $excel.ActiveWorkbook

$statusDone = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns for the b9d5d801 row ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusDone
$overview.Range("C3").Value = $statusDone

# --- zh-cn sheet: Status + Latest Handback DateTime for the b9d5d801 row ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusDone
$zhcn.Range("H3").Value = "2016-03-19 08:39:16"

# --- de-de sheet: Status + Latest Handback DateTime for the b9d5d801 row ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusDone
$dede.Range("H3").Value = "2016-03-19 08:39:21"
